$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin name / link cells (several ranks shifted as the rankings refreshed)
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

# Update price / volume(1h) cells. Force text format first so the numeric-looking
# strings (e.g. "312.70", "1.08%") keep their exact original formatting/precision
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '312.70'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.08%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '37.94'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.09%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.139'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.08%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07922'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.50%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.924'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-2.71%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.248'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.54%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.820'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-7.10%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9281'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.36%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1205'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-7.39%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1933'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.85%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09266'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '5.10%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03348'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-2.23%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09641'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.00%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001370'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.27%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005923'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.41%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.542'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.27%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.407'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.99%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.85%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.291'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '5.71%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1285'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.85%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2589'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4.12%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.04378'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.46%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001253'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '2.83%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004286'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-7.14%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001300'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-3.80%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02115'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-7.23%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05104'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.47%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007618'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.55%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009128'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-7.37%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1358'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.02%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002050'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '0.39%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.008691'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.11%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006690'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '1.51%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.09%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002890'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-3.66%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001200'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.04%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.09%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.09%'
